$d = $word.ActiveDocument

# --- Step 1: locate the run of text that needs to be split + wrapped in
#     proofErr spell-check markers ("asdffgg" -> "A" + "sdffgg") ---------
$r1 = $d.Content
$found = $r1.Find.Execute("asdffgg", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # $r1 now spans exactly the found text ("asdffgg"), courtesy of Find.Execute
    $xml1 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>A</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>sdffgg</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $r1.InsertXML($xml1)
}

# --- Step 2: add a new paragraph after it containing "42refsxfbv " ------
$p1 = $d.Paragraphs.Item(1)
$endOfP1 = $p1.Range.End
$r2 = $d.Range($endOfP1, $endOfP1)
$r2.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = "42refsxfbv "
